$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("evr")

# Update existing values (columns B and C) for rows 2-6
$ws.Range("B2").Value = 0.19089907854464311
$ws.Range("C2").Value = 0.7711273666050199

$ws.Range("B3").Value = 0.1548028361363773

$ws.Range("B4").Value = 0.10973750927077019

$ws.Range("B5").Value = 0.090832384580718489

$ws.Range("B6").Value = 0.078671193901490466

# Add two new rows of data
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 0.074679562968796978

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 0.071504801202223356

# Apply the same style (style index 1) used by other data cells to the new rows
$ws.Range("A6:B6").Copy() | Out-Null
$ws.Range("A7:B7").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:B6").Copy() | Out-Null
$ws.Range("A8:B8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
